# Update "想去人数" (column F) values across sheets to reflect newly
# generated output numbers (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7268
$ws1.Range("F7").Value = 4714
$ws1.Range("F8").Value = 6952
$ws1.Range("F10").Value = 261
$ws1.Range("F21").Value = 1128
$ws1.Range("F23").Value = 45
$ws1.Range("F28").Value = 40
$ws1.Range("F29").Value = 155

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 724
$ws3.Range("F8").Value = 1507
$ws3.Range("F9").Value = 2377

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 724
$ws4.Range("F10").Value = 7268
$ws4.Range("F12").Value = 4714
$ws4.Range("F14").Value = 6952
$ws4.Range("F15").Value = 261
$ws4.Range("F23").Value = 1128
$ws4.Range("F26").Value = 45
